# Update the "Resultado (%)" column (C2:C5) on the Governança sheet with
# new indicator results, and move the active cell selection to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Governança")

# Update the result values used both by the table and the chart (which
# references Governança!$C$2:$C$5 directly, so its cache refreshes too).
$ws.Range("C2").Value = 75
$ws.Range("C3").Value = 20
$ws.Range("C4").Value = 30
$ws.Range("C5").Value = 65

# Force the embedded chart to refresh its cached values from the new data
# (the chart's series reads straight from Governança!$C$2:$C$5).
$excel.CalculateFull()
foreach ($co in $ws.ChartObjects()) {
    $co.Chart.Refresh()
}

# Make sure the sheet is the active one, then move/select the new active cell.
$ws.Activate()
$ws.Range("G3").Select()
